# Update CDA Logical model for ST.r2b
# - Bump Version and Date values on the Metadata sheet
# - Insert a new "Jurisdiction" property row (with an empty value) right after "Contact"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# 1. Update the "Version" value (row 3, column B)
$ws.Cells.Item(3, 2).Value = "2.0.1-sd-202510-matchbox-patch"

# 2. Update the "Date" value (row 8, column B)
$ws.Cells.Item(8, 2).Value = "2025-10-29T22:15:57+01:00"

# 3. Insert a new row after "Contact" (row 10) for the "Jurisdiction" property
$ws.Rows.Item(11).Insert()

# Copy formatting from the row above (Contact) so the new row matches the
# existing style used throughout the property table.
$ws.Cells.Item(10, 1).Copy()
$ws.Cells.Item(11, 1).PasteSpecial(-4122)
$ws.Cells.Item(10, 2).Copy()
$ws.Cells.Item(11, 2).PasteSpecial(-4122)

$ws.Cells.Item(11, 1).Value = "Jurisdiction"
$ws.Cells.Item(11, 2).Value = ""
